$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2293.0
$ws.Range("G2").Value = 2111.0
$ws.Range("H2").Value = 7120.0
$ws.Range("O2").Value = 251.0
$ws.Range("Q2").Value = 29124.0
$ws.Range("D6").Value = 4093.0
$ws.Range("F6").Value = 1266.0
$ws.Range("G6").Value = 1292.0
$ws.Range("H6").Value = 3293.0
$ws.Range("I6").Value = 507.0
$ws.Range("O6").Value = 306.0
$ws.Range("P6").Value = 472.0
$ws.Range("Q6").Value = 17203.0
$ws.Range("F7").Value = 2991.0
$ws.Range("G7").Value = 2673.0
$ws.Range("H7").Value = 7083.0
$ws.Range("P7").Value = 120.0
$ws.Range("D8").Value = 4482.0
$ws.Range("F8").Value = 1949.0
$ws.Range("G8").Value = 1833.0
$ws.Range("H8").Value = 7737.0
$ws.Range("I8").Value = 1126.0
$ws.Range("O8").Value = 304.0
$ws.Range("P8").Value = 280.0
$ws.Range("Q8").Value = 25688.0
$ws.Range("D9").Value = 4468.0
$ws.Range("F9").Value = 2523.0
$ws.Range("G9").Value = 2345.0
$ws.Range("H9").Value = 11022.0
$ws.Range("I9").Value = 2428.0
$ws.Range("J9").Value = 1925.0
$ws.Range("O9").Value = 67.0
$ws.Range("P9").Value = 480.0
$ws.Range("Q9").Value = 3423.0
$ws.Range("D10").Value = 4537.0
$ws.Range("F10").Value = 2274.0
$ws.Range("G10").Value = 1763.0
$ws.Range("H10").Value = 9613.0
$ws.Range("O10").Value = 924.0
$ws.Range("Q10").Value = 92221.0
$ws.Range("F11").Value = 1402.0
$ws.Range("G11").Value = 1235.0
$ws.Range("H11").Value = 5710.0
$ws.Range("I11").Value = 789.0
$ws.Range("J11").Value = 1238.0
$ws.Range("O11").Value = 677.0
$ws.Range("P11").Value = 400.0
$ws.Range("Q11").Value = 39325.0
$ws.Range("M12").Value = 26.0
$ws.Range("P12").Value = 520.0
$ws.Range("D14").Value = 4239.0
$ws.Range("F14").Value = 2317.0
$ws.Range("G14").Value = 2109.0
$ws.Range("H14").Value = 6793.0
$ws.Range("I14").Value = 1091.0
$ws.Range("P14").Value = 310.0
$ws.Range("H15").Value = 3677.0
$ws.Range("O15").Value = 345.0
$ws.Range("P15").Value = 360.0
$ws.Range("Q15").Value = 31542.0
$ws.Range("F17").Value = 2843.0
$ws.Range("G17").Value = 2596.0
$ws.Range("H17").Value = 6744.0
$ws.Range("O20").Value = 635.0
$ws.Range("P20").Value = 560.0
$ws.Range("Q20").Value = 69250.0
$ws.Range("D22").Value = 3975.0
$ws.Range("G22").Value = 1421.0
$ws.Range("H22").Value = 4642.0
$ws.Range("T22").Value = "Arena 12"
$ws.Range("H24").Value = 18708.0
$ws.Range("M24").Value = 16173.0
$ws.Range("P24").Value = 520.0
$ws.Range("D25").Value = 3908.0
$ws.Range("F25").Value = 1382.0
$ws.Range("G25").Value = 1398.0
$ws.Range("H25").Value = 4000.0
$ws.Range("P25").Value = 80.0
$ws.Range("G26").Value = 3521.0
$ws.Range("H26").Value = 12441.0
$ws.Range("J26").Value = 2776.0
$ws.Range("P26").Value = 360.0
$ws.Range("F27").Value = 3746.0
$ws.Range("H27").Value = 11473.0
$ws.Range("P27").Value = 120.0
$ws.Range("F28").Value = 865.0
$ws.Range("G28").Value = 709.0
$ws.Range("H28").Value = 2600.0
$ws.Range("I28").Value = 351.0
$ws.Range("O28").Value = 142.0
$ws.Range("P28").Value = 96.0
$ws.Range("H30").Value = 6751.0
$ws.Range("I30").Value = 1476.0
$ws.Range("O30").Value = 532.0
$ws.Range("P30").Value = 480.0
$ws.Range("Q30").Value = 14014.0
$ws.Range("D31").Value = 4135.0
$ws.Range("F31").Value = 697.0
$ws.Range("H31").Value = 3472.0
$ws.Range("P31").Value = 280.0
$ws.Range("F32").Value = 2171.0
$ws.Range("H32").Value = 7046.0
$ws.Range("I32").Value = 1587.0
$ws.Range("F33").Value = 5230.0
$ws.Range("G33").Value = 5323.0
$ws.Range("H33").Value = 12504.0
$ws.Range("I33").Value = 2274.0
$ws.Range("J33").Value = 1764.0
$ws.Range("O33").Value = 72.0
$ws.Range("Q33").Value = 59389.0
$ws.Range("F34").Value = 964.0
$ws.Range("H34").Value = 4276.0
$ws.Range("I34").Value = 1005.0
$ws.Range("P34").Value = 190.0
$ws.Range("F35").Value = 594.0
$ws.Range("G35").Value = 520.0
$ws.Range("H35").Value = 3518.0
$ws.Range("O35").Value = 317.0
$ws.Range("P35").Value = 480.0
$ws.Range("Q35").Value = 17780.0
$ws.Range("D36").Value = 3961.0
$ws.Range("F36").Value = 3579.0
$ws.Range("G36").Value = 4147.0
$ws.Range("H36").Value = 9857.0
$ws.Range("O36").Value = 319.0
$ws.Range("P36").Value = 360.0
$ws.Range("Q36").Value = 15144.0
$ws.Range("T36").Value = "Arena 12"
$ws.Range("F37").Value = 338.0
$ws.Range("G37").Value = 194.0
$ws.Range("H37").Value = 679.0
$ws.Range("I37").Value = 128.0
$ws.Range("J37").Value = 666.0
$ws.Range("O37").Value = 264.0
$ws.Range("P37").Value = 320.0
$ws.Range("Q37").Value = 3077.0
$ws.Range("D39").Value = 3943.0
$ws.Range("F39").Value = 658.0
$ws.Range("G39").Value = 376.0
$ws.Range("H39").Value = 1757.0
$ws.Range("O39").Value = 96.0
$ws.Range("P39").Value = 120.0
$ws.Range("Q39").Value = 8109.0
$ws.Range("D40").Value = 4329.0
$ws.Range("F40").Value = 1106.0
$ws.Range("H40").Value = 3166.0
$ws.Range("I40").Value = 665.0
$ws.Range("O40").Value = 206.0
$ws.Range("P40").Value = 320.0
$ws.Range("Q40").Value = 20572.0
$ws.Range("T40").Value = "League 2"
$ws.Range("H41").Value = 27112.0
$ws.Range("I41").Value = 3693.0
$ws.Range("P41").Value = 240.0
$ws.Range("D42").Value = 4168.0
$ws.Range("F42").Value = 2441.0
$ws.Range("G42").Value = 2366.0
$ws.Range("H42").Value = 6160.0
$ws.Range("I42").Value = 1197.0
$ws.Range("O42").Value = 316.0
$ws.Range("P42").Value = 280.0
$ws.Range("Q42").Value = 22578.0
$ws.Range("D43").Value = 3879.0
$ws.Range("F43").Value = 1301.0
$ws.Range("G43").Value = 1133.0
$ws.Range("H43").Value = 4392.0
$ws.Range("O43").Value = 778.0
$ws.Range("P43").Value = 480.0
$ws.Range("Q43").Value = 29241.0
$ws.Range("H44").Value = 7847.0
$ws.Range("O44").Value = 166.0
$ws.Range("Q44").Value = 32014.0
$ws.Range("O45").Value = 297.0
$ws.Range("P45").Value = 230.0
$ws.Range("Q45").Value = 22960.0
$ws.Range("F48").Value = 11154.0
$ws.Range("G48").Value = 11855.0
$ws.Range("H48").Value = 26570.0
$ws.Range("O48").Value = 1153.0
$ws.Range("P48").Value = 520.0
$ws.Range("Q48").Value = 77858.0
$ws.Range("H49").Value = 15731.0
$ws.Range("O49").Value = 413.0
$ws.Range("Q49").Value = 33940.0
